$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are plain text (inline strings), e.g. "335.98" / "2.05%".
# Excel auto-converts numeric-/percent-looking text typed into a cell, so each
# touched cell is pre-formatted as Text ("@") to keep the literal string value.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "335.98"
$ws.Range("E2").Value = "2.05%"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "44.02"
$ws.Range("E3").Value = "6.71%"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.803"
$ws.Range("E4").Value = "2.95%"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08343"
$ws.Range("E5").Value = "2.19%"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "8.830"
$ws.Range("E6").Value = "0.93%"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.985"
$ws.Range("E7").Value = "-1.69%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.68%"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9405"
$ws.Range("E9").Value = "2.16%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.39%"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1965"
$ws.Range("E11").Value = "0.34%"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09703"
$ws.Range("E12").Value = "3.44%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "20.56%"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1069"
$ws.Range("E14").Value = "0.90%"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001306"
$ws.Range("E15").Value = "0.47%"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005988"
$ws.Range("E16").Value = "-2.43%"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.497"
$ws.Range("E17").Value = "1.52%"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "4.507"
$ws.Range("E18").Value = "0.23%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.73%"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "8.775"
$ws.Range("E20").Value = "5.56%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.78%"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04400"
$ws.Range("E23").Value = "0.12%"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").Value = "0.52%"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004394"
$ws.Range("E25").Value = "1.26%"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001262"
$ws.Range("E26").Value = "4.98%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003992"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05724"
$ws.Range("E40").Value = "6.04%"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007952"
$ws.Range("E41").Value = "6.67%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.00%"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008957"
$ws.Range("E43").Value = "0.04%"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").Value = "-3.22%"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01057"
$ws.Range("E45").Value = "-8.14%"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007285"
$ws.Range("E46").Value = "10.57%"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003242"
$ws.Range("E48").Value = "1.27%"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").Value = "-0.45%"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.05%"
